$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a temporary text number format on the target cells so that Excel does not
# auto-convert numeric-looking strings (e.g. "0.07980" or "243.43") into actual
# numbers (which would drop trailing zeros / use scientific notation). After the
# values are written, restore each cell's style to match its row's default style
# (taken from column B, which is never touched) so no stray style index remains
# on the edited cells.
$targetCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "E24",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D38",
    "E38",
    "D39",
    "E39",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($cellAddr in $targetCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.453.88'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.884.53'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '0.7202'
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").Value = '243.43'
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.07980'
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").Value = '0.3158'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("D10").Value = '25.01'
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("D11").Value = '0.08167'
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("D12").Value = '1.917.55'
$ws.Range("E12").Value = '  +2.11%  '
$ws.Range("D13").Value = '94.95'
$ws.Range("E13").Value = '  +4.22%  '
$ws.Range("D14").Value = '5.254'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '0.7130'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '6.411'
$ws.Range("E16").Value = '  +4.96%  '
$ws.Range("D17").Value = '0.000008456'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").Value = '29.449.99'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '252.09'
$ws.Range("E19").Value = '  +4.88%  '
$ws.Range("D20").Value = '13.35'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '2.137.85'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '7.792'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = '9.091'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '162.54'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '18.91'
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("D29").Value = '1.508'
$ws.Range("D30").Value = '4.438'
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").Value = '4.298'
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").Value = '1.227'
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("D33").Value = '0.05339'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("D34").Value = '1.951'
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").Value = '0.7565'
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("D36").Value = '1.184'
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("D38").Value = '0.01888'
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").Value = '1.278.16'
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").Value = '6.487'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '113.24'
$ws.Range("E42").Value = '  +3.92%  '
$ws.Range("D43").Value = '75.13'
$ws.Range("E43").Value = '  +3.87%  '
$ws.Range("D44").Value = '0.9097'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").Value = '0.00000000132'
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '2.031.32'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").Value = '1.811'
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = '0.5199'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '9.543'
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("D51").Value = '0.4379'
$ws.Range("E51").Value = '  +0.93%  '

# Restore the default (unstyled) appearance for each touched cell using its row's
# column-B cell as the style reference.
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Style = $ws.Range("B2").Style
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Style = $ws.Range("B3").Style
$ws.Range("E4").Style = $ws.Range("B4").Style
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Style = $ws.Range("B5").Style
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Style = $ws.Range("B6").Style
$ws.Range("E7").Style = $ws.Range("B7").Style
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Style = $ws.Range("B8").Style
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Style = $ws.Range("B9").Style
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Style = $ws.Range("B10").Style
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Style = $ws.Range("B11").Style
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Style = $ws.Range("B12").Style
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Style = $ws.Range("B13").Style
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Style = $ws.Range("B14").Style
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Style = $ws.Range("B15").Style
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Style = $ws.Range("B16").Style
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Style = $ws.Range("B17").Style
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Style = $ws.Range("B18").Style
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Style = $ws.Range("B19").Style
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Style = $ws.Range("B20").Style
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Style = $ws.Range("B21").Style
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Style = $ws.Range("B22").Style
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Style = $ws.Range("B23").Style
$ws.Range("E24").Style = $ws.Range("B24").Style
$ws.Range("E25").Style = $ws.Range("B25").Style
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Style = $ws.Range("B26").Style
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Style = $ws.Range("B27").Style
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Style = $ws.Range("B28").Style
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Style = $ws.Range("B30").Style
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Style = $ws.Range("B31").Style
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Style = $ws.Range("B32").Style
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Style = $ws.Range("B33").Style
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Style = $ws.Range("B34").Style
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Style = $ws.Range("B35").Style
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Style = $ws.Range("B36").Style
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Style = $ws.Range("B38").Style
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Style = $ws.Range("B39").Style
$ws.Range("E40").Style = $ws.Range("B40").Style
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Style = $ws.Range("B41").Style
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Style = $ws.Range("B42").Style
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Style = $ws.Range("B43").Style
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Style = $ws.Range("B44").Style
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Style = $ws.Range("B45").Style
$ws.Range("E46").Style = $ws.Range("B46").Style
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Style = $ws.Range("B47").Style
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Style = $ws.Range("B48").Style
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Style = $ws.Range("B49").Style
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Style = $ws.Range("B50").Style
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Style = $ws.Range("B51").Style
